$p = $ppt.ActivePresentation

# The table with the problem list lives on slide 6 ("표 4" / Table 4).
$s = $p.Slides.Item(6)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cell = $tbl.Cell($r, $c)
                $tr = $cell.Shape.TextFrame.TextRange
                if ($tr.Text -eq "19221197") {
                    $tr.Text = "1922"
                }
            }
        }
    }
}
